$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.893344666666666
$ws.Range("H2").Value = 5.680033999999999
$ws.Range("I2").Value = 0.05525983881677096
$ws.Range("J2").Value = 0.05525983881677096
$ws.Range("M2").Value = 0.007957000000000001
$ws.Range("N2").Value = 0.023871
$ws.Range("O2").Value = 0.0002448939493579708
$ws.Range("P2").Value = 0.0002448939493579708
$ws.Range("Q2").Value = 0.01506534351266667
$ws.Range("R2").Value = 0.135588091614
$ws.Range("S2").Value = [double]"1.353280016872393E-05"
$ws.Range("T2").Value = [double]"1.353280016872393E-05"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.893344666666666
$ws.Range("H3").Value = 5.680033999999999
$ws.Range("I3").Value = 0.05525983881677096
$ws.Range("J3").Value = 0.05525983881677096
$ws.Range("O3").Value = 0.003249135679578298
$ws.Range("P3").Value = 0.003249135679578299
$ws.Range("Q3").Value = 0.1998797653451111
$ws.Range("R3").Value = 1.798917888106
$ws.Range("S3").Value = 0.0001795467139473163
$ws.Range("T3").Value = 0.0001795467139473164
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.893344666666666
$ws.Range("H4").Value = 5.680033999999999
$ws.Range("I4").Value = 0.05525983881677096
$ws.Range("J4").Value = 0.05525983881677096
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02449766666666667
$ws.Range("N4").Value = 0.073493
$ws.Range("O4").Value = 0.0007539688752111494
$ws.Range("P4").Value = 0.0007539688752111494
$ws.Range("Q4").Value = 0.04638252652911111
$ws.Range("R4").Value = 0.4174427387619999
$ws.Range("S4").Value = [double]"4.166419851703021E-05"
$ws.Range("T4").Value = [double]"4.166419851703021E-05"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.893344666666666
$ws.Range("H5").Value = 5.680033999999999
$ws.Range("I5").Value = 0.05525983881677096
$ws.Range("J5").Value = 0.05525983881677096
$ws.Range("M5").Value = 32.353591
$ws.Range("N5").Value = 97.060773
$ws.Range("O5").Value = 0.9957520014958525
$ws.Range("P5").Value = 0.9957520014958525
$ws.Range("Q5").Value = 61.25649896736466
$ws.Range("R5").Value = 551.3084907062819
$ws.Range("S5").Value = 0.05502509510413788
$ws.Range("T5").Value = 0.05502509510413788
$ws.Range("G6").Value = 4.159773333333334
$ws.Range("I6").Value = 0.1214086415227279
$ws.Range("J6").Value = 0.1214086415227279
$ws.Range("M6").Value = 0.007957000000000001
$ws.Range("N6").Value = 0.023871
$ws.Range("O6").Value = 0.0002448939493579708
$ws.Range("P6").Value = 0.0002448939493579708
$ws.Range("Q6").Value = 0.03309931641333334
$ws.Range("R6").Value = 0.29789384772
$ws.Range("S6").Value = [double]"2.973224170868695E-05"
$ws.Range("T6").Value = [double]"2.973224170868695E-05"
$ws.Range("G7").Value = 4.159773333333334
$ws.Range("I7").Value = 0.1214086415227279
$ws.Range("J7").Value = 0.1214086415227279
$ws.Range("O7").Value = 0.003249135679578298
$ws.Range("P7").Value = 0.003249135679578299
$ws.Range("R7").Value = 3.952312957880001
$ws.Range("S7").Value = 0.0003944731489806265
$ws.Range("T7").Value = 0.0003944731489806266
$ws.Range("G8").Value = 4.159773333333334
$ws.Range("I8").Value = 0.1214086415227279
$ws.Range("J8").Value = 0.1214086415227279
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02449766666666667
$ws.Range("N8").Value = 0.073493
$ws.Range("O8").Value = 0.0007539688752111494
$ws.Range("P8").Value = 0.0007539688752111494
$ws.Range("Q8").Value = 0.1019047405288889
$ws.Range("R8").Value = 0.9171426647600002
$ws.Range("S8").Value = [double]"9.153833688980479E-05"
$ws.Range("T8").Value = [double]"9.153833688980479E-05"
$ws.Range("G9").Value = 4.159773333333334
$ws.Range("I9").Value = 0.1214086415227279
$ws.Range("J9").Value = 0.1214086415227279
$ws.Range("M9").Value = 32.353591
$ws.Range("N9").Value = 97.060773
$ws.Range("O9").Value = 0.9957520014958525
$ws.Range("P9").Value = 0.9957520014958525
$ws.Range("Q9").Value = 134.5836050793733
$ws.Range("R9").Value = 1211.25244571436
$ws.Range("S9").Value = 0.1208928977951488
$ws.Range("T9").Value = 0.1208928977951488
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.240212333333333
$ws.Range("H10").Value = 3.720637
$ws.Range("I10").Value = 0.03619728348733726
$ws.Range("J10").Value = 0.03619728348733727
$ws.Range("M10").Value = 0.007957000000000001
$ws.Range("N10").Value = 0.023871
$ws.Range("O10").Value = 0.0002448939493579708
$ws.Range("P10").Value = 0.0002448939493579708
$ws.Range("Q10").Value = 0.009868369536333334
$ws.Range("R10").Value = 0.088815325827
$ws.Range("S10").Value = [double]"8.864495709244083E-06"
$ws.Range("T10").Value = [double]"8.864495709244084E-06"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.240212333333333
$ws.Range("H11").Value = 3.720637
$ws.Range("I11").Value = 0.03619728348733726
$ws.Range("J11").Value = 0.03619728348733727
$ws.Range("O11").Value = 0.003249135679578298
$ws.Range("P11").Value = 0.003249135679578299
$ws.Range("Q11").Value = 0.1309288026258889
$ws.Range("R11").Value = 1.178359223633
$ws.Range("S11").Value = 0.0001176098852825179
$ws.Range("T11").Value = 0.0001176098852825179
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.240212333333333
$ws.Range("H12").Value = 3.720637
$ws.Range("I12").Value = 0.03619728348733726
$ws.Range("J12").Value = 0.03619728348733727
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.02449766666666667
$ws.Range("N12").Value = 0.073493
$ws.Range("O12").Value = 0.0007539688752111494
$ws.Range("P12").Value = 0.0007539688752111494
$ws.Range("Q12").Value = 0.03038230833788889
$ws.Range("R12").Value = 0.273440775041
$ws.Range("S12").Value = [double]"2.729162511664679E-05"
$ws.Range("T12").Value = [double]"2.729162511664679E-05"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.240212333333333
$ws.Range("H13").Value = 3.720637
$ws.Range("I13").Value = 0.03619728348733726
$ws.Range("J13").Value = 0.03619728348733727
$ws.Range("M13").Value = 32.353591
$ws.Range("N13").Value = 97.060773
$ws.Range("O13").Value = 0.9957520014958525
$ws.Range("P13").Value = 0.9957520014958525
$ws.Range("Q13").Value = 40.12532258582233
$ws.Range("R13").Value = 361.127903272401
$ws.Range("S13").Value = 0.03604351748122885
$ws.Range("T13").Value = 0.03604351748122886
$ws.Range("G14").Value = 26.96925
$ws.Range("H14").Value = 80.90774999999999
$ws.Range("I14").Value = 0.7871342361731639
$ws.Range("J14").Value = 0.7871342361731638
$ws.Range("M14").Value = 0.007957000000000001
$ws.Range("N14").Value = 0.023871
$ws.Range("O14").Value = 0.0002448939493579708
$ws.Range("P14").Value = 0.0002448939493579708
$ws.Range("Q14").Value = 0.21459432225
$ws.Range("R14").Value = 1.93134890025
$ws.Range("S14").Value = 0.0001927644117713158
$ws.Range("T14").Value = 0.0001927644117713158
$ws.Range("G15").Value = 26.96925
$ws.Range("H15").Value = 80.90774999999999
$ws.Range("I15").Value = 0.7871342361731639
$ws.Range("J15").Value = 0.7871342361731638
$ws.Range("O15").Value = 0.003249135679578298
$ws.Range("P15").Value = 0.003249135679578299
$ws.Range("Q15").Value = 2.84713473275
$ws.Range("R15").Value = 25.62421259475
$ws.Range("S15").Value = 0.002557505931367838
$ws.Range("T15").Value = 0.002557505931367838
$ws.Range("G16").Value = 26.96925
$ws.Range("H16").Value = 80.90774999999999
$ws.Range("I16").Value = 0.7871342361731639
$ws.Range("J16").Value = 0.7871342361731638
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.02449766666666667
$ws.Range("N16").Value = 0.073493
$ws.Range("O16").Value = 0.0007539688752111494
$ws.Range("P16").Value = 0.0007539688752111494
$ws.Range("Q16").Value = 0.66068369675
$ws.Range("R16").Value = 5.94615327075
$ws.Range("S16").Value = 0.0005934747146876676
$ws.Range("T16").Value = 0.0005934747146876675
$ws.Range("G17").Value = 26.96925
$ws.Range("H17").Value = 80.90774999999999
$ws.Range("I17").Value = 0.7871342361731639
$ws.Range("J17").Value = 0.7871342361731638
$ws.Range("M17").Value = 32.353591
$ws.Range("N17").Value = 97.060773
$ws.Range("O17").Value = 0.9957520014958525
$ws.Range("P17").Value = 0.9957520014958525
$ws.Range("Q17").Value = 872.55208407675
$ws.Range("R17").Value = 7852.968756690749
$ws.Range("S17").Value = 0.7837904911153371
$ws.Range("T17").Value = 0.7837904911153369
